# V19: Killing ghosts give no reward
# Adds a new row (24) to Sheet1 with the results for the v19-18300 run.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("A24").Value = "v19-18300"
$ws.Range("B24").Value = 8071
$ws.Range("C24").Value = "Not including ghosts in reward"
$ws.Range("D24").Value = 166
$ws.Range("E24").Value = 0.66
$ws.Range("F24").Value = 1.78
$ws.Range("G24").Value = 0
$ws.Range("H24").Value = 5
$ws.Range("I24").Value = 3061.81
$ws.Range("J24").Value = 2000
$ws.Range("K24").Value = 6400
$ws.Range("L24").Value = 238.52
$ws.Range("M24").Value = 192
$ws.Range("N24").Value = 244
$ws.Range("O24").Value = 2.04
$ws.Range("P24").Value = 0
$ws.Range("Q24").Value = 7
$ws.Range("R24").Value = 107.5
$ws.Range("S24").Value = 71.3
$ws.Range("T24").Value = 188.6

# Match the number formats used by the rest of the table columns.
$ws.Range("D24").NumberFormat = "0"
$ws.Range("E24").NumberFormat = "0.00%"
$ws.Range("F24").NumberFormat = "0.00"
$ws.Range("G24").NumberFormat = "0"
$ws.Range("H24").NumberFormat = "0"
$ws.Range("I24").NumberFormat = "0.00"
$ws.Range("J24").NumberFormat = "0"
$ws.Range("K24").NumberFormat = "0"
$ws.Range("L24").NumberFormat = "0.00"
$ws.Range("M24").NumberFormat = "0"
$ws.Range("N24").NumberFormat = "0"
$ws.Range("O24").NumberFormat = "0.00"
$ws.Range("P24").NumberFormat = "0"
$ws.Range("Q24").NumberFormat = "0"
$ws.Range("R24").NumberFormat = "0.00"
$ws.Range("S24").NumberFormat = "0.00"
$ws.Range("T24").NumberFormat = "0.00"

# Reflect the selection that was active when the workbook was saved.
$ws.Range("C23").Select()
